$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Shikuku Emmanuel"
$ws.Range("B8").Value = "Nabwana"
$ws.Range("C8").Value = "Kenyan"
$ws.Range("D8").Value = "Marketing"
$ws.Range("E8").Value = "Marketing"
